$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update component positions / rotations (1.1 revision: R14 & R15 VUSB
#     sense CP2102, plus the neighbouring C21, R3 and R11 placements that
#     moved as part of the same re-layout) ---

# C21 (row 22)
$ws.Range("B22").Value = 126.43000000000001
$ws.Range("C22").Value = -88.609999999999999
$ws.Range("D22").Value = 180

# R3 (row 38)
$ws.Range("B38").Value = 135
$ws.Range("C38").Value = -91.629999999999995

# R11 (row 46)
$ws.Range("B46").Value = 125.70999999999999
$ws.Range("C46").Value = -91.629999999999995
$ws.Range("D46").Value = 180

# R14 (row 49)
$ws.Range("B49").Value = 130.27000000000001
$ws.Range("C49").Value = -104.94

# R15 (row 50)
$ws.Range("B50").Value = 130.27000000000001
$ws.Range("C50").Value = -103.43000000000001
$ws.Range("D50").Value = 0

# --- Header row is no longer bold; the numeric header cells (Mid X / Mid Y /
#     Rotation) pick up the same "0.000000" number format already used by the
#     data columns beneath them ---
$ws.Range("A1:E1").Font.Bold = $false
$ws.Range("A1:E1").Font.Name = "Sans"
$ws.Range("A1:E1").Font.Size = 10
$ws.Range("A1:E1").Font.Color = 0
$ws.Range("B1:D1").NumberFormat = "0.000000"

# Keep the data range's formatting explicit too (regular weight, same numeric
# format) so it stays consistent with the header after the restyle.
$ws.Range("B2:D55").Font.Bold = $false
$ws.Range("B2:D55").NumberFormat = "0.000000"

# --- Column A reverts to the sheet's default (auto) width ---
$ws.Columns.Item(1).ColumnWidth = 8.1
